# Update countries & provincias Spain
# - Refresh the "last updated" timestamp.
# - Fix country label ordering (Belgica/Panama, Bermudas/Brunei/Islas Turcas y
#   Caicos, Santa Lucia/Timor Oriental) together with each country's updated
#   statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp refresh (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 05:23"

# Kazajistan (row 29) - updated stats, same country
$ws.Range("B29").Value = 97829
$ws.Range("C29").Value = 907
$ws.Range("D29").Value = 71609
$ws.Range("E29").Value = 25162

# Rows 39/40: Panama <-> Belgica swap position, each with its own new stats
$ws.Range("A39").Value = "Belgica"
$ws.Range("B39").Value = 72784
$ws.Range("C39").Value = 768
$ws.Range("D39").Value = 17728
$ws.Range("E39").Value = 45190
$ws.Range("G39").Value = 5
$ws.Range("H39").Value = 9866

$ws.Range("A40").Value = "Panama"
$ws.Range("B40").Value = 72560
$ws.Range("D40").Value = 46675
$ws.Range("E40").Value = 24294
$ws.Range("H40").Value = 1591

# Honduras (row 50) - updated stats, same country
$ws.Range("B50").Value = 46365
$ws.Range("C50").Value = 610
$ws.Range("D50").Value = 6355
$ws.Range("E50").Value = 38545
$ws.Range("G50").Value = 19
$ws.Range("H50").Value = 1465

# Jamaica (row 153) - updated stats, same country
$ws.Range("B153").Value = 987
$ws.Range("C153").Value = 29
$ws.Range("E153").Value = 229
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 13

# Rows 184/185/186: Bermudas -> Islas Turcas y Caicos -> Brunei cascade,
# each landing row carries its predecessor's old stats except row 184 which
# gets genuinely new "Islas Turcas y Caicos" numbers.
$ws.Range("A184").Value = "Islas Turcas y Caicos"
$ws.Range("B184").Value = 170
$ws.Range("C184").Value = 29
$ws.Range("D184").Value = 39
$ws.Range("E184").Value = 129
$ws.Range("H184").Value = 2

$ws.Range("A185").Value = "Bermudas"
$ws.Range("B185").Value = 157
$ws.Range("D185").Value = 144
$ws.Range("E185").Value = 4
$ws.Range("H185").Value = 9

$ws.Range("A186").Value = "Brunei"
$ws.Range("B186").Value = 142
$ws.Range("D186").Value = 138
$ws.Range("E186").Value = 1
$ws.Range("H186").Value = 3

# Rows 202/203: Timor Oriental <-> Santa Lucia swap (identical stats, label
# order only).
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"
